$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.352.65'
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").Value = '3.448.57'
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '408.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.78'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +13.44%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.610'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.30%  '
$ws.Range("D8").Value = '3.440.17'
$ws.Range("E8").Value = '  +0.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.684'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.130'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +7.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.97'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.34%  '
$ws.Range("E13").Value = '  -1.46%  '
$ws.Range("D14").Value = '3.987.43'
$ws.Range("E14").Value = '  +0.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.71'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.82%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.02'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.71%  '
$ws.Range("D17").Value = '3.481.18'
$ws.Range("E17").Value = '  +0.16%  '
$ws.Range("D18").Value = '62.347.21'
$ws.Range("E18").Value = '  +0.21%  '
$ws.Range("E19").Value = '  -2.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.22'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.93%  '
$ws.Range("E21").Value = '  +11.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.28'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.34%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '84.67'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.13'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '310.71'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("E26").Value = '  -1.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.66'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '30.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.53'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.59'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.69%  '
$ws.Range("E31").Value = '  -1.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.117'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '43.79'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +10.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.75'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.82%  '
$ws.Range("E35").Value = '  +1.87%  '
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0489'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.69%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.36'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.997'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("E40").Value = '  -0.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.02'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.127'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '137.41'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.98'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.21%  '
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.288'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.08%  '
$ws.Range("B46").Value = 'Celestia'
$ws.Range("C46").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.07'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.93'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.60%  '
$ws.Range("E48").Value = '  -2.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '21.81'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.63%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.164.17'
$ws.Range("E50").Value = '  -0.84%  '
$ws.Range("B51").Value = 'ApeXProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.37'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.54%  '
